# Removed texture img files and display
# - Insert a new finding row (row 14) about texture images hurting text
#   readability, pushing the existing rows 14-19 down to 15-20.
# - Update the "Défaut de conception du CSS" row (now row 20) wording to
#   be specific to the contact-form CSS, and tweak its supporting text.
# - Update the "Fichier CSS trop grand" row (now row 19) wording and fill
#   in the previously empty E column with a placeholder/whitespace cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row 14 (pushes old rows 14-19 down to 15-20) ---
$ws.Rows(14).Insert()
$ws.Rows(14).RowHeight = 67.2

$ws.Range("A14").Value = "Accessibilité"
$ws.Range("B14").Value = "Les img de textures gènent la lecture des textes"
$ws.Range("C14").Value = "Les textures trop prononcées peuvent gommer les contrastes entre le texte et l'arrière plan, ce qui peut entrainer des gènes de lectures auprès des utilisateurs"
$ws.Range("D14").Value = "Ne pas utiliser d'images de textures ou alors un grain très léger et avec peu d'opacité"
$ws.Range("E14").Value = "Retirer les textures existantes"
$ws.Range("G14").Value = "OK"

# --- Row 19 (previously row 18): "Fichier CSS trop grand" ---
$ws.Range("C19").Value = "Implique un temps supplémentaire qu'un dev mettra à rechercher les segments à modifier pour mettre à jour ou maintenir le site"
$ws.Range("E19").Value = "                                               "
$ws.Rows(19).RowHeight = 68.4

# --- Row 20 (previously row 19): "Défaut de conception du CSS" ---
$ws.Range("B20").Value = "Défaut de conception du formulaire de contact (CSS)"
$ws.Range("D20").Value = "Pour une meilleure maintenabilité du site web, il est préférable d'organiser certains éléments dans une flexbox contenue dans une grid"

# The "OK" cell for the newly inserted row ends up one row below where the
# selection previously was (G14 -> G15).
[void]$ws.Range("G15").Select()

Write-Output "edit applied"
